$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A60").Value = "TestTitle"
$ws.Range("A60").Font.Bold = $true
Write-Host "done"
